$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(237671823369, 'MFS ENTREE COLLEGE MALANGUE', 'Rte_1', 'Cite Sic Stade Marion Ocm', 189785.8, 103700, -86085.79999999999, 0.5464054739606441, 'Cite Sic')
    ,@(237672128028, 'DELVIN NDIFON BAH', 'Rte_1', 'Cite Sic Stade Marion Ocm', 25030, 64935, 39905, 2.594286855773072, 'Cite Sic')
    ,@(237672277367, 'TOP MOBIL KM5 LTDLA_POLAS_BTQ_KM5', 'Rte_0', 'Cite Sic Stade Marion Ocm', 26250, 1100121, 1073871, 41.90937142857143, 'Cite Sic')
    ,@(237674853971, 'NJOSSEU TCHOUNZOU TOP MOBILE', 'Rte_1', 'Cite Sic Stade Marion Ocm', 182169.7, 379307, 197137.3, 2.082162950260114, 'Cite Sic')
    ,@(237674884705, 'BAH AMADOU MOUNTAGHA ETS MOBILE FINANCIAL SERVICES MFS', 'Rte_0', 'Cite Sic Stade Marion Ocm', 5000, 19074, 14074, 3.8148, 'Cite Sic')
    ,@(237675779272, 'RODES NGWEM KEMAYOU', 'Rte_1', 'Cite Sic Stade Marion Ocm', 184162.0370370371, 137655, -46507.03703703705, 0.7474667538148269, 'Cite Sic')
    ,@(237677304210, 'CARINE OROCK', 'Rte_2', 'Cite Sic Stade Marion Ocm', 61193.33333333334, 90031, 28837.66666666666, 1.471255038675237, 'Cite Sic')
    ,@(237678267353, 'LA NEGRESSE SARL EMBOLA BELTUS MBU', 'Rte_0', 'Cite Sic Stade Marion Ocm', 16520, 0, -16520, 0, 'Cite Sic')
    ,@(237678370615, 'ESSEN ONGOLONG BERTHE HORTENSE ETS MOBILE FINANCIAL SERVICES MFS', 'Rte_1', 'Cite Sic Stade Marion Ocm', 100229.8, 123243, 23013.2, 1.229604369159671, 'Cite Sic')
    ,@(237678836319, 'KAMDOM DOMINIQUE STEPHANIE ETS MOBILE FINANCIAL SERVICES MFS', 'Rte_1', 'Cite Sic Stade Marion Ocm', 25619.09090909091, 109522, 83902.90909090909, 4.275015081083, 'Cite Sic')
    ,@(237678922502, 'NWOAGA TCHAMDJOU EPSE KAMSEU EMILINE ETS LE CONTENT', 'Rte_0', 'Cite Sic Stade Marion Ocm', 56220, 563253, 507033, 10.01872998932764, 'Cite Sic')
    ,@(237679884264, 'MFS CICAM', 'Rte_1', 'Cite Sic Stade Marion Ocm', 59610, 117365, 57755, 1.968881060224795, 'Cite Sic')
    ,@(237681019523, 'ETS MOULAY RIPERT AND COMPANY', 'Rte_1', 'Cite Sic Stade Marion Ocm', 57805.71428571429, 90676, 32870.28571428571, 1.568633847370502, 'Cite Sic')
    ,@(237681125655, 'EMILE MADELO', 'Rte_1', 'Cite Sic Stade Marion Ocm', 47813.75, 12821, -34992.75, 0.2681446236699694, 'Cite Sic')
    ,@(237681240793, 'MBANE EMILIE FRANCOISE ETS MOBILE FINANCIAL SERVICES MFS', 'Rte_0', 'Cite Sic Stade Marion Ocm', 7438.090909090909, 2959, -4479.090909090909, 0.3978171329397817, 'Cite Sic')
    ,@(237682117915, 'MEKUEKO FOUDJO BERLINE DIDIANE ETS MOBILE FINANCIAL SERVICES MFS', 'Rte_1', 'Cite Sic Stade Marion Ocm', 123740, 23687, -100053, 0.1914255697430095, 'Cite Sic')
    ,@(237682154553, 'N A ISUFUH MIEMONA NGESSY ETS MOBILE FINANCIAL SERVICES MFS', 'Rte_2', 'Cite Sic Stade Marion Ocm', 136800, 31551, -105249, 0.2306359649122807, 'Cite Sic')
    ,@(237682803277, 'NGUEMASSOM RENE MARTIAL LA NEGRESSE SARL', 'Rte_0', 'Cite Sic Stade Marion Ocm', 59869.80000000002, 69625, 9755.199999999983, 1.162940247002662, 'Cite Sic')
    ,@(237683323481, 'ETS LE CONTENT TSAZE DONFOUET FLORETTE ROSINE', 'Rte_0', 'Cite Sic Stade Marion Ocm', 25924, 509243, 483319, 19.64368924548681, 'Cite Sic')
    ,@(237683368985, 'MFS BELL HENRIE BERNARD', 'Rte_1', 'Cite Sic Stade Marion Ocm', 51230, 232132, 180902, 4.531173140737849, 'Cite Sic')
    ,@(237683432110, 'ERODINE TOUMENI', 'Rte_0', 'Cite Sic Stade Marion Ocm', 15600, 822255, 806655, 52.70865384615384, 'Cite Sic')
    ,@(237683743490, 'ETS LE CONTENT NGAH MARIE', 'Rte_0', 'Cite Sic Stade Marion Ocm', 7776.875, 4823, -2953.875, 0.6201719842481717, 'Cite Sic')
)

$startRow = 86
for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $r = $startRow + $idx
    $row = $data[$idx]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}
